# "Optimized Loop Code LDB once"
# Swap the Op Code/Operand rows for "OUT" (row 2) and "LDB" (row 3) on the
# "Loop Code" sheet, bump two operand flags from 0 to 1, and move the
# selection to C14.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Loop Code")

# Row 2 becomes LDB / 1111 (was OUT / 1110)
$ws.Range("B2").Value = "LDB"
$ws.Range("C2").Value = 1111

# Row 3 becomes OUT / 1110 (was LDB / 1111)
$ws.Range("B3").Value = "OUT"
$ws.Range("C3").Value = 1110

# Operand flags flipped from 0 to 1
$ws.Range("C7").Value = 1
$ws.Range("C11").Value = 1

# Update the active selection
$ws.Activate()
$ws.Range("C14").Select()
